$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap H1 / I1 labels ---
$ws.Range("H1").Value = "high_detail"
$ws.Range("I1").Value = "n_ratings"

$rotation = "[0.9961567660501535, 0.008693328396189522, 0.0007605676829745856, 0.08715242412403446]"
$position = "[0.58, 0.0, 0.65]"

# --- Row 2: extended ---
$ws.Range("B2").Value = $rotation
$ws.Range("C2").Value = "splats/mcmc-truck-extended-1.ksplat"
$ws.Range("D2").Value = "splats/default-truck-extended-1.ksplat"
$ws.Range("E2").Value = "extended"
$ws.Range("F2").Value = "truck"
$ws.Range("G2").Value = $position
$ws.Range("H2").Value = $true
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = "5dwojHXzFd2RzJ2ZPNtD"

# --- Row 3: medium ---
$ws.Range("B3").Value = $rotation
$ws.Range("C3").Value = "splats/mcmc-truck-medium-1.ksplat"
$ws.Range("D3").Value = "splats/default-truck-medium-1.ksplat"
$ws.Range("E3").Value = "medium"
$ws.Range("F3").Value = "truck"
$ws.Range("G3").Value = $position
$ws.Range("H3").Value = $false
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = "BGOjoywll3sSUYo2pBwY"

# --- Row 4: low ---
$ws.Range("B4").Value = $rotation
$ws.Range("C4").Value = "splats/mcmc-truck-low-1.ksplat"
$ws.Range("D4").Value = "splats/default-truck-low-1.ksplat"
$ws.Range("E4").Value = "low"
$ws.Range("F4").Value = "truck"
$ws.Range("G4").Value = $position
$ws.Range("H4").Value = $false
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = "OyhrkDZsWy64SkqfZAY5"

# --- Row 5: high ---
$ws.Range("B5").Value = $rotation
$ws.Range("C5").Value = "splats/mcmc-truck-high-1.ksplat"
$ws.Range("D5").Value = "splats/default-truck-high-1.ksplat"
$ws.Range("E5").Value = "high"
$ws.Range("F5").Value = "truck"
$ws.Range("G5").Value = $position
$ws.Range("H5").Value = $true
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = "PqpflZDmB5tBiM15v9MQ"
